$d = $word.ActiveDocument

$replacements = @(
    @("684×8=", "423×7="),
    @("149×2=", "101×6="),
    @("848×5=", "113×2="),
    @("824×8=", "312×5="),
    @("217×7=", "720×7="),
    @("488×7=", "839×8="),
    @("847×8=", "374×4="),
    @("901×7=", "854×4="),
    @("775×9=", "607×2="),
    @("462×7=", "267×3="),
    @("588×3=", "827×8="),
    @("238×6=", "882×2="),
    @("514×3=", "600×3="),
    @("745×5=", "551×2="),
    @("650×8=", "638×7="),
    @("816×8=", "561×8="),
    @("812×2=", "999×8="),
    @("609×2=", "925×9="),
    @("259×2=", "243×5="),
    @("746×2=", "613×8="),
    @("129×3=", "395×3="),
    @("872×7=", "307×3="),
    @("954×4=", "392×4="),
    @("921×6=", "663×2="),
    @("803×8=", "631×8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
